$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, update the last existing row (currently row 6) in place: ---
# B6/C6 get new "LUBY" input files, A6 renamed TestCase9 -> TestCase6, D6 -> No
$ws.Range("B6").Value = "D:\Options_10001933_LUBY_3_2000Jun30_updated.pdf"
$ws.Range("C6").Value = "D:\Options_10001933_LUBY_3_2000Jun30.pdf"
$ws.Range("A6").Value = "TestCase6"
$ws.Range("D6").Value = "No"

# Also flip the ExecutionMode for rows 4 and 5 (TestCase4 / TestCase5) to No
$ws.Range("D4").Value = "No"
$ws.Range("D5").Value = "No"

# --- Now insert a new row above current row 2, shifting rows 2-6 down to 3-7 ---
$ws.Rows.Item(2).Insert()

# Row 2 (new): TestCase1
$ws.Range("A2").Value = "TestCase1"
$ws.Range("B2").Value = "D:\Benchmark 1–French.pdf"
$ws.Range("C2").Value = "D:\Benchmark 1–French_New.pdf"
$ws.Range("D2").Value = "Yes"
$ws.Range("A2:D2").WrapText = $true
$ws.Range("B2:C2").ClearFormats()

# Row 3 (was row 2): TestCase2, with swapped files
$ws.Range("A3").Value = "TestCase2"
$ws.Range("B3").Value = "D:\Benchmark 1–French_New.pdf"
$ws.Range("C3").Value = "D:\Benchmark 1–French.pdf"
$ws.Range("D3").Value = "Yes"
$ws.Range("A3:D3").WrapText = $true
$ws.Range("B3:C3").ClearFormats()

# Row 4 (was row 3): TestCase3 - unchanged content, row height becomes 30
$ws.Range("A4:D4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 30

# Row 5 (was row 4): TestCase4
$ws.Range("A5:D5").WrapText = $true

# Row 6 (was row 5): TestCase5
$ws.Range("A6:D6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 30

# Row 7 (was row 6): TestCase6
$ws.Range("A7:D7").WrapText = $true
$ws.Rows.Item(7).AutoFit()

$ws.Range("D4").Select()
